# Insert one new data row right before the current row 354. This shifts the
# existing rows 354-423 down to 355-424 (Excel's normal Insert behaviour),
# expanding the used range from A1:T423 to A1:T424, and fills the freshly
# inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(354).Insert()

$ws.Cells.Item(354, 1).Value = 11
$ws.Cells.Item(354, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(354, 3).Value = 'Bíobío'
$ws.Cells.Item(354, 4).Value = 44637
$ws.Cells.Item(354, 5).Value = 8
$ws.Cells.Item(354, 6).Value = 'Fruta'
$ws.Cells.Item(354, 7).Value = 100102
$ws.Cells.Item(354, 8).Value = 'Cítricos'
$ws.Cells.Item(354, 9).Value = 100102003
$ws.Cells.Item(354, 10).Value = 'Limón'
$ws.Cells.Item(354, 11).Value = 'Sin especificar'
$ws.Cells.Item(354, 12).Value = '1a plateado'
$ws.Cells.Item(354, 13).Value = 270
$ws.Cells.Item(354, 14).Value = 22000
$ws.Cells.Item(354, 15).Value = 23000
$ws.Cells.Item(354, 16).Value = 22556
$ws.Cells.Item(354, 17).Value = '$/malla 16 kilos'
$ws.Cells.Item(354, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(354, 19).Value = 1410
$ws.Cells.Item(354, 20).Value = 16
